# component design.pptx edit
#   - flexbox/device-label declarations updated to include all breakpoints
#   - breadcrumbs wireframe (old design) removed from slide 2
#   - "PRODUCT CARD COMPONENT" comments renamed to "FEATURE COMPONENT"
#   - cached "today" date field re-stamped on master + all layouts

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2: breadcrumbs component slide
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# "Tablet, Desktop" -> "Mobile, Tablet, Desktop" (device-support label)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $sh = $s2.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "Tablet, Desktop") {
            $sh.TextFrame.TextRange.Text = "Mobile, Tablet, Desktop"
        }
    }
}

# Remove the old breadcrumbs wireframe shapes (superseded by new designs)
$idsToRemove = @(32, 33, 34, 38, 39, 40, 41, 44, 45, 46, 47)
$shapesToRemove = @()
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $sh = $s2.Shapes.Item($i)
    if ($idsToRemove -contains $sh.Id) {
        $shapesToRemove += $sh
    }
}
foreach ($sh in $shapesToRemove) {
    $sh.Delete()
}

# ---------------------------------------------------------------------------
# Slide 4: product card component slide
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
for ($i = 1; $i -le $s4.Shapes.Count; $i++) {
    $sh = $s4.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $t = $sh.TextFrame.TextRange.Text
        if ($t -eq "Tablet, Desktop (Mobile on next slide)") {
            $sh.TextFrame.TextRange.Text = "Desktop (Mobile and Tablet on next slide)"
        } elseif ($t -eq "/* PRODUCT CARD COMPONENT */") {
            $sh.TextFrame.TextRange.Text = "/* FEATURE COMPONENT */"
        }
    }
}

# ---------------------------------------------------------------------------
# Slide 5: product card component (mobile) slide
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
for ($i = 1; $i -le $s5.Shapes.Count; $i++) {
    $sh = $s5.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $t = $sh.TextFrame.TextRange.Text
        if ($t -eq "Mobile") {
            $sh.TextFrame.TextRange.Text = "Mobile, Tablet"
        } elseif ($t -eq "/* PRODUCT CARD COMPONENT */") {
            $sh.TextFrame.TextRange.Text = "/* FEATURE COMPONENT */"
        }
    }
}

# ---------------------------------------------------------------------------
# Re-stamp the cached "datetimeFigureOut" date placeholder text on the
# slide master and every slide layout (9/30/21 -> 10/18/21)
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapeCollection) {
    for ($i = 1; $i -le $shapeCollection.Count; $i++) {
        $sh = $shapeCollection.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "9/30/21") {
                $sh.TextFrame.TextRange.Text = "10/18/21"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
